# Add two new columns to the sheet: I ("I0") and J ("IF"), matching the
# layout/commit message "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold / border / centered) of the existing header
# cell H1 onto the new header cells I1:J1, so the style matches the rest
# of row 1 exactly (reuses the same style definition instead of creating
# a new one).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-41: row number, I value, J value.
$data = @(
    @(2, 7, 8),
    @(3, 7, 7),
    @(4, 4, 4),
    @(5, 7, 7),
    @(6, 6, 6),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 7, 7),
    @(11, 6, 6),
    @(12, 7, 8),
    @(13, 4, 4),
    @(14, 6, 6),
    @(15, 8, 8),
    @(16, 6, 6),
    @(17, 8, 8),
    @(18, 6, 6),
    @(19, 9, 9),
    @(20, 6, 6),
    @(21, 7, 7),
    @(22, 7, 7),
    @(23, 7, 7),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 6, 7),
    @(27, 6, 6),
    @(28, 8, 9),
    @(29, 4, 4),
    @(30, 8, 8),
    @(31, 7, 7),
    @(32, 9, 9),
    @(33, 7, 7),
    @(34, 9, 9),
    @(35, 6, 6),
    @(36, 7, 8),
    @(37, 2, 2),
    @(38, 7, 7),
    @(39, 9, 9),
    @(40, 5, 5),
    @(41, 6, 7)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
